# Refresh cryptocurrency price/volume figures (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.811.05'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.275.88'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.641'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '78.64'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.643'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.11'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.71%  '
$ws.Range('D14').Value = '2.615.23'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.870'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = '2.277.74'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '42.699.44'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '0.0₃0996'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '234.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.59%  '
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0856'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.96%  '
$ws.Range('E33').Value = '  -5.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.128'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.57'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0303'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.62'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.26'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.93%  '
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '115.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +17.45%  '
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.44%  '
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.28%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  -3.59%  '
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.92%  '
